$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that look numeric (e.g. "4.30", "0.0000133").
# Excel's native type-inference on Range.Value would coerce these to numbers
# and normalize/round them (losing the exact text, e.g. trailing zeros or
# multi-dot "thousands" formatting). Prefix with a leading apostrophe so Excel
# keeps them as literal text, exactly like typing '4.30 into a cell would.

$ws.Range("D2").Value = "'59.047.05"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "'2.588.14"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'522.86"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "'139.18"
$ws.Range("E6").Value = "  -2.13%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "'0.565"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "'2.600.01"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("E11").Value = "  -0.45%  "
$ws.Range("E12").Value = "  +1.16%  "
$ws.Range("E13").Value = "  +3.27%  "
$ws.Range("D14").Value = "'3.042.94"
$ws.Range("E14").Value = "  +0.35%  "
$ws.Range("D15").Value = "'58.997.17"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'2.603.10"
$ws.Range("E17").Value = "  +2.07%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000133"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").Value = "'338.44"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").Value = "'4.30"
$ws.Range("E20").Value = "  -0.16%  "
$ws.Range("D21").Value = "'10.09"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("D22").Value = "'6.51"
$ws.Range("E22").Value = "  +2.57%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("E27").Value = "  +0.12%  "
$ws.Range("D28").Value = "'7.02"
$ws.Range("E28").Value = "  +0.35%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "'0.0₃0724"
$ws.Range("E30").Value = "  -2.68%  "
$ws.Range("E31").Value = "  -5.15%  "
$ws.Range("E32").Value = "  +0.20%  "
$ws.Range("D33").Value = "'18.70"
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("E36").Value = "  -1.61%  "
$ws.Range("D37").Value = "'36.80"
$ws.Range("E37").Value = "  +2.38%  "
$ws.Range("E38").Value = "  +1.34%  "
$ws.Range("D39").Value = "'0.826"
$ws.Range("E39").Value = "  -1.16%  "
$ws.Range("E40").Value = "  -5.71%  "
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("D43").Value = "'271.71"
$ws.Range("E43").Value = "  +0.66%  "
$ws.Range("E44").Value = "  +0.84%  "
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("E46").Value = "  +0.36%  "
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("D48").Value = "'18.40"
$ws.Range("E48").Value = "  -2.15%  "
$ws.Range("D49").Value = "'1.967.00"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("D50").Value = "'4.51"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("E51").Value = "  -0.32%  "
